$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (F column) counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 253
$ws1.Range("F4").Value = 14
$ws1.Range("F5").Value = 6527
$ws1.Range("F6").Value = 5294
$ws1.Range("F7").Value = 441
$ws1.Range("F8").Value = 65
$ws1.Range("F11").Value = 225

# Sheet "全部类型" - update "想去人数" (F column) counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 253
$ws4.Range("F4").Value = 14
$ws4.Range("F5").Value = 6528
$ws4.Range("F6").Value = 5294
$ws4.Range("F7").Value = 441
$ws4.Range("F8").Value = 65
$ws4.Range("F11").Value = 225
